$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.01591511301927033
$ws.Range("C3").Value = 0.01540387312590868
$ws.Range("D3").Value = 0.01542400039369613

# Row 4 - GradientBoostingRegressor
$ws.Range("B4").Value = 0.01262962350598083
$ws.Range("C4").Value = 0.01266261103080567
$ws.Range("D4").Value = 0.01258595684638744

# Row 5 - AdaBoostRegressor
$ws.Range("B5").Value = 0.04097081456282165
$ws.Range("C5").Value = 0.04190027278907638
$ws.Range("D5").Value = 0.04634824795316145
